$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.88903588170909
$ws.Range("C2").Value = 8.914391069298858
$ws.Range("E2").Value = 14.00576996002193
$ws.Range("F2").Value = 47.13236582543433
$ws.Range("G2").Value = 54.74629978412564
$ws.Range("H2").Value = 20.77934469939177
$ws.Range("J2").Value = 10.02121827833715
$ws.Range("K2").Value = 10.89431346536594
$ws.Range("L2").Value = 11.25414926259406
$ws.Range("M2").Value = 16.17100675595452
$ws.Range("B3").Value = 14.76818376998755
$ws.Range("C3").Value = 8.891888479506246
$ws.Range("E3").Value = 14.02366839923798
$ws.Range("F3").Value = 47.10553734332376
$ws.Range("G3").Value = 54.69275038002143
$ws.Range("H3").Value = 20.80723857698386
$ws.Range("J3").Value = 10.02819323177236
$ws.Range("K3").Value = 10.81310739745977
$ws.Range("L3").Value = 11.26459044306134
$ws.Range("M3").Value = 16.16794370854588
$ws.Range("B4").Value = 14.69718008317833
$ws.Range("C4").Value = 8.877636286502833
$ws.Range("E4").Value = 14.03619889554529
$ws.Range("F4").Value = 47.0976629261466
$ws.Range("G4").Value = 54.67054597351212
$ws.Range("H4").Value = 20.8271399549611
$ws.Range("J4").Value = 10.0326193798214
$ws.Range("K4").Value = 10.76559128574097
$ws.Range("L4").Value = 11.27244862291934
$ws.Range("M4").Value = 16.16869874344243
$ws.Range("B5").Value = 14.6690802622842
$ws.Range("C5").Value = 8.871718996135849
$ws.Range("E5").Value = 14.04169329617436
$ws.Range("F5").Value = 47.09661973717634
$ws.Range("G5").Value = 54.66418517905828
$ws.Range("H5").Value = 20.83594718919653
$ws.Range("J5").Value = 10.03445930531933
$ws.Range("K5").Value = 10.74683718491995
$ws.Range("L5").Value = 11.27601552425845
$ws.Range("M5").Value = 16.16967136250648
$ws.Range("B6").Value = 14.66446557485152
$ws.Range("C6").Value = 8.870729829847573
$ws.Range("E6").Value = 14.04262909816098
$ws.Range("F6").Value = 47.0965774005037
$ws.Range("G6").Value = 54.6632913656825
$ws.Range("H6").Value = 20.83745172350317
$ws.Range("J6").Value = 10.03476701640973
$ws.Range("K6").Value = 10.74376041759236
$ws.Range("L6").Value = 11.27662984779976
$ws.Range("M6").Value = 16.16987308339743
$ws.Range("B7").Value = 14.69679770065733
$ws.Range("C7").Value = 8.877556926686145
$ws.Range("E7").Value = 14.03627142257579
$ws.Range("F7").Value = 47.09764008471918
$ws.Range("G7").Value = 54.67044930431667
$ws.Range("H7").Value = 20.82725590983212
$ws.Range("J7").Value = 10.03264404679669
$ws.Range("K7").Value = 10.76533587020069
$ws.Range("L7").Value = 11.27249525013116
$ws.Range("M7").Value = 16.1687091655767
$ws.Range("B8").Value = 14.84672042671911
$ws.Range("C8").Value = 8.906722627130859
$ws.Range("E8").Value = 14.01162196728649
$ws.Range("F8").Value = 47.12133324676322
$ws.Range("G8").Value = 54.72562283114134
$ws.Range("H8").Value = 20.78838646398836
$ws.Range("J8").Value = 10.02359355062519
$ws.Range("K8").Value = 10.8658389782963
$ws.Range("L8").Value = 11.25744944209444
$ws.Range("M8").Value = 16.16940493229505
$ws.Range("B9").Value = 15.16471714053727
$ws.Range("C9").Value = 8.960457785640068
$ws.Range("E9").Value = 13.97547963504825
$ws.Range("F9").Value = 47.23581886025903
$ws.Range("G9").Value = 54.91830492656671
$ws.Range("H9").Value = 20.73419319943108
$ws.Range("J9").Value = 10.00697699427834
$ws.Range("K9").Value = 11.08062452634117
$ws.Range("L9").Value = 11.23939568951802
$ws.Range("M9").Value = 16.19157887967024
$ws.Range("B10").Value = 15.41103056065296
$ws.Range("C10").Value = 8.997826731026077
$ws.Range("E10").Value = 13.95631910438758
$ws.Range("F10").Value = 47.36106602443358
$ws.Range("G10").Value = 55.11098375201556
$ws.Range("H10").Value = 20.70782629110015
$ws.Range("J10").Value = 9.995448126905991
$ws.Range("K10").Value = 11.2479443022912
$ws.Range("L10").Value = 11.23306718221957
$ws.Range("M10").Value = 16.22040274370795
$ws.Range("B11").Value = 15.52540490571523
$ws.Range("C11").Value = 9.014368346247736
$ws.Range("E11").Value = 13.94919851314428
$ws.Range("F11").Value = 47.4268799609888
$ws.Range("G11").Value = 55.20961671438972
$ws.Range("H11").Value = 20.69875392275077
$ws.Range("J11").Value = 9.990348594518444
$ws.Range("K11").Value = 11.32584370116215
$ws.Range("L11").Value = 11.23168339664038
$ws.Range("M11").Value = 16.23619854965706
$ws.Range("B12").Value = 15.56901262779637
$ws.Range("C12").Value = 9.020566318887811
$ws.Range("E12").Value = 13.94673073529565
$ws.Range("F12").Value = 47.45306256993808
$ws.Range("G12").Value = 55.24853221309274
$ws.Range("H12").Value = 20.69573862388007
$ws.Range("J12").Value = 9.98843822083405
$ws.Range("K12").Value = 11.35557418881199
$ws.Range("L12").Value = 11.23137335396501
$ws.Range("M12").Value = 16.24256218282799
$ws.Range("B13").Value = 15.55960833453677
$ws.Range("C13").Value = 9.019234411532157
$ws.Range("E13").Value = 13.94725206012905
$ws.Range("F13").Value = 47.44736780719798
$ws.Range("G13").Value = 55.24008168656073
$ws.Range("H13").Value = 20.69636933430685
$ws.Range("J13").Value = 9.988848735244845
$ws.Range("K13").Value = 11.34916129298281
$ws.Range("L13").Value = 11.23143062797419
$ws.Range("M13").Value = 16.24117472996147
$ws.Range("B14").Value = 15.52898681105702
$ws.Range("C14").Value = 9.014879580404186
$ws.Range("E14").Value = 13.9489909101678
$ws.Range("F14").Value = 47.42900882710635
$ws.Range("G14").Value = 55.21278703048236
$ws.Range("H14").Value = 20.69849743093203
$ws.Range("J14").Value = 9.990191012715071
$ws.Range("K14").Value = 11.32828514715659
$ws.Range("L14").Value = 11.23165360712612
$ws.Range("M14").Value = 16.23671445444685
$ws.Range("B15").Value = 15.51026775130465
$ws.Range("C15").Value = 9.012203523544567
$ws.Range("E15").Value = 13.95008575635076
$ws.Range("F15").Value = 47.41792721268624
$ws.Range("G15").Value = 55.19627167400579
$ws.Range("H15").Value = 20.69985567407328
$ws.Range("J15").Value = 9.991015889041833
$ws.Range("K15").Value = 11.31552733015788
$ws.Range("L15").Value = 11.23181802158707
$ws.Range("M15").Value = 16.23403205040961
$ws.Range("B16").Value = 15.40359949889005
$ws.Range("C16").Value = 8.996736473638977
$ws.Range("E16").Value = 13.95681649114446
$ws.Range("F16").Value = 47.35694189778703
$ws.Range("G16").Value = 55.10475770221086
$ws.Range("H16").Value = 20.70847799539397
$ws.Range("J16").Value = 9.995784298995439
$ws.Range("K16").Value = 11.24288717980328
$ws.Range("L16").Value = 11.23318761469697
$ws.Range("M16").Value = 16.21942413428675
$ws.Range("B17").Value = 15.33873077366522
$ws.Range("C17").Value = 8.987130595888825
$ws.Range("E17").Value = 13.96135369005143
$ws.Range("F17").Value = 47.32178592195685
$ws.Range("G17").Value = 55.05142038456402
$ws.Range("H17").Value = 20.71451597088157
$ws.Range("J17").Value = 9.998746601777937
$ws.Range("K17").Value = 11.1987642585267
$ws.Range("L17").Value = 11.23441003442459
$ws.Range("M17").Value = 16.21114738824663
$ws.Range("B18").Value = 15.30164113434307
$ws.Range("C18").Value = 8.98156247731273
$ws.Range("E18").Value = 13.96411360607684
$ws.Range("F18").Value = 47.30239776057456
$ws.Range("G18").Value = 55.02177683644008
$ws.Range("H18").Value = 20.71826389923033
$ws.Range("J18").Value = 10.0004640972676
$ws.Range("K18").Value = 11.17355543272481
$ws.Range("L18").Value = 11.23525392771368
$ws.Range("M18").Value = 16.20663966970862
$ws.Range("B19").Value = 15.28912234405429
$ws.Range("C19").Value = 8.979669815301436
$ws.Range("E19").Value = 13.96507389417237
$ws.Range("F19").Value = 47.29597657068194
$ws.Range("G19").Value = 55.01191813924508
$ws.Range("H19").Value = 20.71958012021097
$ws.Range("J19").Value = 10.00104796176156
$ws.Range("K19").Value = 11.16505002720901
$ws.Range("L19").Value = 11.23556386640811
$ws.Range("M19").Value = 16.20515697173825
$ws.Range("B20").Value = 15.34561354575217
$ws.Range("C20").Value = 8.988157616912117
$ws.Range("E20").Value = 13.96085515407988
$ws.Range("F20").Value = 47.32544223505674
$ws.Range("G20").Value = 55.05699124389211
$ws.Range("H20").Value = 20.71384475148614
$ws.Range("J20").Value = 9.998429847145585
$ws.Range("K20").Value = 11.20344385499988
$ws.Range("L20").Value = 11.23426534175991
$ws.Range("M20").Value = 16.21200231996921
$ws.Range("B21").Value = 15.53797334585369
$ws.Range("C21").Value = 9.016160491773851
$ws.Range("E21").Value = 13.94847396932405
$ws.Range("F21").Value = 47.43436719142645
$ws.Range("G21").Value = 55.22076176986288
$ws.Range("H21").Value = 20.69786095335933
$ws.Range("J21").Value = 9.989796192399842
$ws.Range("K21").Value = 11.33441088998817
$ws.Range("L21").Value = 11.23158231423124
$ws.Range("M21").Value = 16.23801420700737
$ws.Range("B22").Value = 15.66540301458289
$ws.Range("C22").Value = 9.034077344434888
$ws.Range("E22").Value = 13.94171440919892
$ws.Range("F22").Value = 47.51289612660585
$ws.Range("G22").Value = 55.33691150803642
$ws.Range("H22").Value = 20.68986385618747
$ws.Range("J22").Value = 9.984274239064254
$ws.Range("K22").Value = 11.42134405790937
$ws.Range("L22").Value = 11.23107544211694
$ws.Range("M22").Value = 16.25723971260074
$ws.Range("B23").Value = 15.59724732073753
$ws.Range("C23").Value = 9.024550030978194
$ws.Range("E23").Value = 13.94520048996397
$ws.Range("F23").Value = 47.47031596641395
$ws.Range("G23").Value = 55.27409114282442
$ws.Range("H23").Value = 20.69390797335724
$ws.Range("J23").Value = 9.987210418235474
$ws.Range("K23").Value = 11.3748320183806
$ws.Range("L23").Value = 11.23123225680436
$ws.Range("M23").Value = 16.24677641905381
$ws.Range("B24").Value = 15.34250120538182
$ws.Range("C24").Value = 8.987693442666064
$ws.Range("E24").Value = 13.96108007048725
$ws.Range("F24").Value = 47.3237866499646
$ws.Range("G24").Value = 55.05446947754953
$ws.Range("H24").Value = 20.71414734808091
$ws.Range("J24").Value = 9.998573006948863
$ws.Range("K24").Value = 11.20132771528021
$ws.Range("L24").Value = 11.23433031770552
$ws.Range("M24").Value = 16.21161502416344
$ws.Range("B25").Value = 15.07631860048388
$ws.Range("C25").Value = 8.946291455683973
$ws.Range("E25").Value = 13.98395565078099
$ws.Range("F25").Value = 47.19759468477811
$ws.Range("G25").Value = 54.85716731601882
$ws.Range("H25").Value = 20.74649314401109
$ws.Range("J25").Value = 10.0113521654925
$ws.Range("K25").Value = 11.02075299007345
$ws.Range("L25").Value = 11.24305822894401
$ws.Range("M25").Value = 16.18336715551591
